$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.3558718861209965
$wsSummary.Range("C2").Value = 0.06510416666666667
$wsSummary.Range("D2").Value = 0.8928571428571429
$wsSummary.Range("E2").Value = 0.1213592233009709
$wsSummary.Range("F2").Value = 0.2520161290322581
$wsSummary.Range("G2").Value = 0.5996309963099631
$wsSummary.Range("H2").Value = 0.8013643659711074
$wsSummary.Range("I2").Value = 25
$wsSummary.Range("J2").Value = 359
$wsSummary.Range("K2").Value = 175
$wsSummary.Range("L2").Value = 3

# --- Sheet: Classification Report ---
$wsClass = $wb.Worksheets.Item("Classification Report")
$wsClass.Range("B2").Value = 0.9831460674157303
$wsClass.Range("C2").Value = 0.3277153558052435
$wsClass.Range("D2").Value = 0.4915730337078651

$wsClass.Range("B3").Value = 0.06510416666666667
$wsClass.Range("C3").Value = 0.8928571428571429
$wsClass.Range("D3").Value = 0.1213592233009709

$wsClass.Range("B4").Value = 0.3558718861209965
$wsClass.Range("C4").Value = 0.3558718861209965
$wsClass.Range("D4").Value = 0.3558718861209965
$wsClass.Range("E4").Value = 0.3558718861209965

$wsClass.Range("B5").Value = 0.5241251170411985
$wsClass.Range("C5").Value = 0.6102862493311931
$wsClass.Range("D5").Value = 0.306466128504418

$wsClass.Range("B6").Value = 0.937407325029656
$wsClass.Range("C6").Value = 0.3558718861209965
$wsClass.Range("D6").Value = 0.4731282175310092

# --- Sheet: Confusion Matrix ---
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 175
$wsConf.Range("C2").Value = 359
$wsConf.Range("B3").Value = 3
$wsConf.Range("C3").Value = 25
